$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2 through 32 in column C had their Fitness value changed to 7293
$ws.Range("C2:C32").Value = 7293
